$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of outage data (row 5) to the Active_Outages sheet.
$ws.Range("A5").Value = ""
$ws.Range("B5").Value = "R4"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "JED0123"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = "SCECO"
$ws.Range("J5").Value = "In progress"
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = "Latis"
